# Weekly data refresh: a new, more recent observation is inserted at row 171
# (pushing the existing rows 171-199 down to 172-200), growing the used range
# from A1:R199 to A1:R200.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 171; existing rows 171..199 shift down to 172..200.
$ws.Rows.Item(171).Insert()

# Populate the newly inserted row 171 with the latest observation.
$ws.Range("A171").Value = 10
$ws.Range("B171").Value = "Vega Modelo de Temuco"
$ws.Range("C171").Value = "La Araucanía"
$ws.Range("D171").Value = 44476
$ws.Range("E171").Value = 9
$ws.Range("F171").Value = 100112009
$ws.Range("G171").Value = "Acelga"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 50
$ws.Range("K171").Value = 8000
$ws.Range("L171").Value = 8000
$ws.Range("M171").Value = 8000
$ws.Range("N171").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O171").Value = "Provincia de Cautín"
$ws.Range("P171").Value = 667
$ws.Range("Q171").Value = 12
$ws.Range("R171").Value = "Hortaliza"
